$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ltp")

# Update LTP (column B) and PREV (column C) values
$ws.Range("B2").Value = 2848.95
$ws.Range("C2").Value = 2809.9

$ws.Range("B3").Value = 454.05
$ws.Range("C3").Value = 438.9

$ws.Range("B4").Value = 1685.8
$ws.Range("C4").Value = 1681.2

$ws.Range("B5").Value = 7327.75
$ws.Range("C5").Value = 7258.65

$ws.Range("B6").Value = 241.4
$ws.Range("C6").Value = 239.35

$ws.Range("B7").Value = 231.1
$ws.Range("C7").Value = 232.95

$ws.Range("B8").Value = 48628.8
$ws.Range("C8").Value = 48822.2

$ws.Range("B9").Value = 726.4
$ws.Range("C9").Value = 728.65

$ws.Range("B10").Value = 4143.5
$ws.Range("C10").Value = 4091.9

$ws.Range("B11").Value = 156.15
$ws.Range("C11").Value = 155.5

$ws.Range("B12").Value = 1466.1
$ws.Range("C12").Value = 1472.45

$ws.Range("B13").Value = 614.85
$ws.Range("C13").Value = 614.3

$ws.Range("B14").Value = 1598.95
$ws.Range("C14").Value = 1610.55

$ws.Range("B15").Value = 748.1
$ws.Range("C15").Value = 736.35

$ws.Range("B16").Value = 536.1
$ws.Range("C16").Value = 527.2

$ws.Range("B17").Value = 1729.4
$ws.Range("C17").Value = 1734.45

$ws.Range("B18").Value = 276.7
$ws.Range("C18").Value = 273.25

$ws.Range("B19").Value = 21885.95
$ws.Range("C19").Value = 21942.3

$ws.Range("B20").Value = 311.15
$ws.Range("C20").Value = 313.9

$ws.Range("B21").Value = 642.05
$ws.Range("C21").Value = 651.4

$ws.Range("B22").Value = 712.2
$ws.Range("C22").Value = 715.55

$ws.Range("B23").Value = 779.95
$ws.Range("C23").Value = 753.9

$ws.Range("B24").Value = 332.15
$ws.Range("C24").Value = 329.8

$ws.Range("B25").Value = 139.6
$ws.Range("C25").Value = 138.15

# Update the active selection on the sheet
$ws.Range("F11").Select()
